$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive (stable) leading text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*随着科学技术的发展*") {
        $target = $p
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found"
} else {
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End

    # The original "_GoBack" bookmark sits (collapsed) right at the start of
    # this paragraph's run content. It is not enumerated by Bookmarks (hidden
    # bookmark) but can still be reached by name. Remove it here since the
    # edit re-creates it at a different position inside the rebuilt runs.
    try {
        $existingGoBack = $d.Bookmarks.Item("_GoBack")
        $existingGoBack.Delete()
    } catch {
        Write-Host "No existing _GoBack bookmark to remove"
    }

    # Exclude the trailing paragraph mark from the replaced range so the
    # paragraph's own pPr (and the paragraph break) stay untouched.
    $r = $d.Range($pStart, $pEnd - 1)

    $fragment = '<w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:hint="eastAsia"/></w:rPr><w:t>随着科学技术的发展，在通信网络建设中</w:t></w:r><w:r><w:t>DWDM</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:hint="eastAsia"/></w:rPr><w:t>技术作为应用最为广泛的一种技术，在应用的过程中，灵活性交叉等缺点和不足需要完善，</w:t></w:r><w:r><w:t>OTN</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:hint="eastAsia"/></w:rPr><w:t>技术作为以</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>DWDM</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:hint="eastAsia"/></w:rPr><w:t>技术为基础的一种技术，在通信网络建设中的应用，需要将</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>DEDM</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:hint="eastAsia"/></w:rPr><w:t>技术与</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>OT</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>N</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="宋体" w:hAnsi="宋体" w:hint="eastAsia"/></w:rPr><w:t>技术有效的结合，提高通信网络建设的质量，实现建设目标。</w:t></w:r>'

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $fragment + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xml)

    Write-Host "Replaced paragraph content"
}
